$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (market day) is inserted at row 21, pushing the
# existing rows 21-36 down to 22-37. Row 21's Insert() call will copy
# formatting (e.g. the date style on column D) from the row above, just
# like Excel's native "Insert" behaviour.
$ws.Rows.Item(21).Insert()

$ws.Cells.Item(21, 1).Value2 = 11
$ws.Cells.Item(21, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value2 = "Bíobío"
$ws.Cells.Item(21, 4).Value2 = 44827
$ws.Cells.Item(21, 5).Value2 = 8
$ws.Cells.Item(21, 6).Value2 = 100112026
$ws.Cells.Item(21, 7).Value2 = "Haba"
$ws.Cells.Item(21, 8).Value2 = "Sin especificar"
$ws.Cells.Item(21, 9).Value2 = "Primera"
$ws.Cells.Item(21, 10).Value2 = 100
$ws.Cells.Item(21, 11).Value2 = 9000
$ws.Cells.Item(21, 12).Value2 = 10000
$ws.Cells.Item(21, 13).Value2 = 9500
$ws.Cells.Item(21, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(21, 16).Value2 = 380
$ws.Cells.Item(21, 17).Value2 = 25
$ws.Cells.Item(21, 18).Value2 = "Hortaliza"
